$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data: B2 changes from "institution_1" to "suisse"
$ws.Range("B2").Value = "suisse"

# Add new row 3 with data
$ws.Range("A3").Value = "pop_com_1"
$ws.Range("B3").Value = "ofs-pop"

# Resize the table to include the new row
$tbl = $ws.ListObjects.Item("Tableau1")
$tbl.Resize($ws.Range("A1:B3"))

# Update selection to reflect where the user last clicked
$ws.Range("B4").Select()
